# Generate Report for Handback
# Row 7 on both the "zh-cn" and "de-de" sheets describes the handback for
# bc8b34db-f44c-4425-870a-178ca68f5aa7.md. A handback file was produced, but
# it was not built from the very latest handoff commit, so the report now
# records: the (outdated) target file that came back, its handback
# timestamp, a hyperlink to the actual commit that was handed back, and an
# error message describing the mismatch between that commit and the latest
# one.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a033aa2c3ad38facc2c7f08901452df9bce92e/e2e/bc8b34db-f44c-4425-870a-178ca68f5aa7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a8574ef76cb7e1ca7b4d38d862d9455439fbd73/e2e/bc8b34db-f44c-4425-870a-178ca68f5aa7.md."
$targetFileLink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a033aa2c3ad38facc2c7f08901452df9bce92e/e2e/bc8b34db-f44c-4425-870a-178ca68f5aa7.md"
$targetFileDisplay = "bc8b34db-f44c-4425-870a-178ca68f5aa7.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I7").Value = $targetFileDisplay
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $targetFileLink, "", "", $targetFileDisplay)

$wsZhCn.Range("J7").Value = "bc8b34db-f44c-4425-870a-178ca68f5aa7.30ec4f7c1b76516dfd9a86c6a89dd9e0a18b2d7c.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-25 18:44:43"
$wsZhCn.Range("P7").Value = $errorDetail

$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I7").Value = $targetFileDisplay
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $targetFileLink, "", "", $targetFileDisplay)

$wsDeDe.Range("J7").Value = "bc8b34db-f44c-4425-870a-178ca68f5aa7.30ec4f7c1b76516dfd9a86c6a89dd9e0a18b2d7c.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-25 18:44:49"
$wsDeDe.Range("P7").Value = $errorDetail

$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
